$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 84 and row 85 (columns A:D)
foreach ($col in @("A", "B", "C", "D")) {
    $addr84 = "$col" + "84"
    $addr85 = "$col" + "85"
    $tmp = $ws.Range($addr84).Value2
    $ws.Range($addr84).Value2 = $ws.Range($addr85).Value2
    $ws.Range($addr85).Value2 = $tmp
}

# Update the active selection to D85
$ws.Range("D85").Select()

# Update workbook view window position
$excel.ActiveWindow.Left = 31755
$excel.ActiveWindow.Top = 2010
